$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.321.59'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.528.61'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.93'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.17'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.74%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.531.44'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0993'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.39'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.975.44'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.263.10'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000141'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.538.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.98'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.18'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.81'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.51'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.419'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.52%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.83'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0771'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.01'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.45'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.50'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.96'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '287.30'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.805'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.997'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.84'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '124.24'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0926'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.61'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0223'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.31%  '
